$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
# Row 12
$ws.Range("H12").Value = 2957
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# Row 43
$ws.Range("H43").Value = 1666.6666
$ws.Range("I43").Value = 1666.6666
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1666.6666
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1597.6666
$ws.Range("N43").ClearContents()
# Row 86
$ws.Range("H86").Value = 4438.25
$ws.Range("I86").Value = 3917.6667
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 3917.6667
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -2794.6667
$ws.Range("N86").Value = -8246
# Row 89
$ws.Range("H89").Value = 4438.25
$ws.Range("I89").Value = 3917.6667
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 19588.3335
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -13972.3335
$ws.Range("N89").Value = -41232
# Row 127
$ws.Range("H127").Value = 3541.6667
$ws.Range("I127").Value = 2083.3333
$ws.Range("K127").Value = 6249.999899999999
$ws.Range("M127").Value = -1289.999899999999
# Row 131
$ws.Range("H131").Value = 1399.6
$ws.Range("I131").Value = 1399.6
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 4198.799999999999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 841.2000000000007
$ws.Range("N131").ClearContents()
# Row 141
$ws.Range("H141").Value = 2773.5
$ws.Range("I141").Value = 2312.5715
$ws.Range("K141").Value = 6937.7145
$ws.Range("M141").Value = -1757.7145

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
# Row 96
$ws.Range("H96").Value = 39999.5
$ws.Range("J96").Value = 39999.5
$ws.Range("L96").Value = 39999.5
$ws.Range("N96").Value = -45491.5
# Row 122
$ws.Range("H122").Value = 2134.2222
$ws.Range("I122").Value = 2261.0667
$ws.Range("K122").Value = 6783.2001
$ws.Range("M122").Value = -4333.2001

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
# Row 64
$ws.Range("H64").Value = 1012.9
$ws.Range("I64").Value = 1148.2
$ws.Range("J64").Value = 877.6
$ws.Range("K64").Value = 1148.2
$ws.Range("L64").Value = 877.6
$ws.Range("M64").Value = -923.2
$ws.Range("N64").Value = -1327.6
# Row 67
$ws.Range("H67").Value = 1012.9
$ws.Range("I67").Value = 1148.2
$ws.Range("J67").Value = 877.6
$ws.Range("K67").Value = 1148.2
$ws.Range("L67").Value = 877.6
$ws.Range("M67").Value = -368.2
$ws.Range("N67").Value = -2437.6
# Row 134
$ws.Range("H134").Value = 6404.1816
$ws.Range("I134").Value = 5294.6
$ws.Range("J134").Value = 17500
$ws.Range("K134").Value = 15883.8
$ws.Range("L134").Value = 52500
$ws.Range("M134").Value = -13348.8
$ws.Range("N134").Value = -57570

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
# Row 16
$ws.Range("H16").Value = 946
$ws.Range("I16").Value = 900.1739
$ws.Range("K16").Value = 900.1739
$ws.Range("M16").Value = -613.1739
# Row 28
$ws.Range("H28").Value = 42666.168
$ws.Range("J28").Value = 42666.168
$ws.Range("L28").Value = 42666.168
$ws.Range("N28").Value = -43156.168
# Row 31
$ws.Range("H31").Value = 8226.727999999999
$ws.Range("I31").Value = 8613.571
$ws.Range("J31").Value = 7549.75
$ws.Range("K31").Value = 8613.571
$ws.Range("L31").Value = 7549.75
$ws.Range("M31").Value = -8318.571
$ws.Range("N31").Value = -8139.75
# Row 34
$ws.Range("H34").Value = 8226.727999999999
$ws.Range("I34").Value = 8613.571
$ws.Range("J34").Value = 7549.75
$ws.Range("K34").Value = 8613.571
$ws.Range("L34").Value = 7549.75
$ws.Range("M34").Value = -8411.571
$ws.Range("N34").Value = -7953.75
# Row 58
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
# Row 113
$ws.Range("H113").Value = 946
$ws.Range("I113").Value = 900.1739
$ws.Range("K113").Value = 900.1739
$ws.Range("M113").Value = 1269.8261
# Row 132
$ws.Range("H132").Value = 2394.1428
$ws.Range("I132").Value = 2129
$ws.Range("K132").Value = 6387
$ws.Range("M132").Value = -3857
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
# Row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
# Row 33
$ws.Range("H33").Value = 569.8461
$ws.Range("I33").Value = 407.83334
$ws.Range("J33").Value = 708.7143
$ws.Range("K33").Value = 2447.00004
$ws.Range("L33").Value = 4252.2858
$ws.Range("M33").Value = -2164.00004
$ws.Range("N33").Value = -4818.2858
# Row 63
$ws.Range("H63").Value = 4386.4
$ws.Range("I63").Value = 3006
$ws.Range("J63").Value = 5306.6665
$ws.Range("K63").Value = 9018
$ws.Range("L63").Value = 15919.9995
$ws.Range("M63").Value = -8269
$ws.Range("N63").Value = -17417.9995
# Row 66
$ws.Range("H66").Value = 4386.4
$ws.Range("I66").Value = 3006
$ws.Range("J66").Value = 5306.6665
$ws.Range("K66").Value = 27054
$ws.Range("L66").Value = 47759.9985
$ws.Range("M66").Value = -23310
$ws.Range("N66").Value = -55247.9985
# Row 98
$ws.Range("H98").Value = 2183.111
$ws.Range("I98").Value = 2482.4
$ws.Range("J98").Value = 1809
$ws.Range("K98").Value = 7447.200000000001
$ws.Range("L98").Value = 5427
$ws.Range("M98").Value = -5949.200000000001
$ws.Range("N98").Value = -8423
# Row 131
$ws.Range("H131").Value = 2272.8667
$ws.Range("I131").Value = 1780.8889
$ws.Range("K131").Value = 5342.6667
$ws.Range("M131").Value = -302.6666999999998
# Row 137
$ws.Range("H137").Value = 5819
$ws.Range("J137").Value = 4705.5
$ws.Range("L137").Value = 14116.5
$ws.Range("N137").Value = -24316.5

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
# Row 3
$ws.Range("H3").Value = 3186791.5
$ws.Range("J3").Value = 1255850.9
$ws.Range("L3").Value = 1255850.9
$ws.Range("N3").Value = -1256082.9
# Row 80
$ws.Range("H80").Value = 4824.7
$ws.Range("I80").Value = 4207.8335
$ws.Range("J80").Value = 5750
$ws.Range("K80").Value = 4207.8335
$ws.Range("L80").Value = 5750
$ws.Range("M80").Value = -3209.8335
$ws.Range("N80").Value = -7746
# Row 83
$ws.Range("H83").Value = 4824.7
$ws.Range("I83").Value = 4207.8335
$ws.Range("J83").Value = 5750
$ws.Range("K83").Value = 21039.1675
$ws.Range("L83").Value = 28750
$ws.Range("M83").Value = -16047.1675
$ws.Range("N83").Value = -38734
# Row 93
$ws.Range("H93").Value = 80000
$ws.Range("J93").Value = 80000
$ws.Range("L93").Value = 80000
$ws.Range("N93").Value = -83744
# Row 97
$ws.Range("H97").Value = 2098.3333
$ws.Range("I97").Value = 2098.3333
$ws.Range("K97").Value = 2098.3333
$ws.Range("M97").Value = -1602.3333
# Row 134
$ws.Range("H134").Value = 39999
$ws.Range("J134").Value = 39999
$ws.Range("L134").Value = 119997
$ws.Range("N134").Value = -125067
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 1585.6
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1732
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1732
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -2322
# Row 27
$ws.Range("H27").Value = 1585.6
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1732
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1732
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1946
# Row 33
$ws.Range("H33").Value = 9008.5
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
# Row 46
$ws.Range("H46").Value = 401999.4
$ws.Range("I46").Value = 2000000
$ws.Range("K46").Value = 2000000
$ws.Range("M46").Value = -1999812
# Row 56
$ws.Range("H56").Value = 45874.75
$ws.Range("J56").Value = 45833.332
$ws.Range("L56").Value = 45833.332
$ws.Range("N56").Value = -47215.332
# Row 101
$ws.Range("H101").Value = 40362
$ws.Range("J101").Value = 40362
$ws.Range("L101").Value = 40362
$ws.Range("N101").Value = -46852
# Row 104
$ws.Range("H104").Value = 28499
$ws.Range("J104").Value = 28499
$ws.Range("L104").Value = 28499
$ws.Range("N104").Value = -35487

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
# Row 92
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
# Row 97
$ws.Range("H97").Value = 61000
$ws.Range("J97").Value = 61000
$ws.Range("L97").Value = 61000
$ws.Range("N97").Value = -62982
# Row 107
$ws.Range("H107").Value = 500
$ws.Range("I107").Value = 500
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 1500
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 420
$ws.Range("N107").Value = -5340
# Row 137
$ws.Range("H137").Value = 74999.336
$ws.Range("J137").Value = 74999.336
$ws.Range("L137").Value = 74999.336
$ws.Range("N137").Value = -85199.336
